$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.203383214053351
$ws.Range("C2").Value = 4923.6
$ws.Range("D2").Value = 0.01535458685751464
$ws.Range("E2").Value = 23.6
$ws.Range("F2").Value = 126
$ws.Range("I2").Value = 0.08197787898503578
$ws.Range("J2").Value = 3.4
$ws.Range("K2").Value = 0.0022121014964216

$ws.Range("B3").Value = 1.775319622012229
$ws.Range("C3").Value = 3193.8
$ws.Range("D3").Value = 0.01634241245136187
$ws.Range("E3").Value = 29.4
$ws.Range("F3").Value = 41
$ws.Range("I3").Value = 0.02279043913285158
$ws.Range("J3").Value = 4.4
$ws.Range("K3").Value = 0.002445803224013341

$ws.Range("B4").Value = 4.58653314427508
$ws.Range("C4").Value = 8017.8
$ws.Range("D4").Value = 0.03247285569866215
$ws.Range("E4").Value = 56.2
$ws.Range("F4").Value = 187.4
$ws.Range("I4").Value = 0.1066721808657292
$ws.Range("J4").Value = 5.4
$ws.Range("K4").Value = 0.003094384707287933

$ws.Range("B5").Value = 7.352278870698461
$ws.Range("C5").Value = 13452.8
$ws.Range("D5").Value = 0.03998785955724499
$ws.Range("E5").Value = 73.40000000000001
$ws.Range("F5").Value = 290.2
$ws.Range("I5").Value = 0.1626965738779845
$ws.Range("J5").Value = 14.6
$ws.Range("K5").Value = 0.008076131122476515
